$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header text in A1 and C1 (no format change needed, both already
#     share the same bold+border style) ---
$a1 = $ws.Cells.Item(1, 1).Value2
$c1 = $ws.Cells.Item(1, 3).Value2
$ws.Cells.Item(1, 1).Value2 = $c1
$ws.Cells.Item(1, 3).Value2 = $a1

# --- Swap column A and column C data + number formats for rows 2-8 ---
$aVals = @{}
$cVals = @{}
$aFmt = @{}
$cFmt = @{}
for ($r = 2; $r -le 8; $r++) {
    $aVals[$r] = $ws.Cells.Item($r, 1).Value2
    $cVals[$r] = $ws.Cells.Item($r, 3).Value2
    $aFmt[$r] = $ws.Cells.Item($r, 1).NumberFormat
    $cFmt[$r] = $ws.Cells.Item($r, 3).NumberFormat
}

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $cVals[$r]
    if ($cFmt[$r] -eq "General") {
        $ws.Cells.Item($r, 1).ClearFormats()
    } else {
        $ws.Cells.Item($r, 1).NumberFormat = $cFmt[$r]
    }
    $ws.Cells.Item($r, 3).Value2 = $aVals[$r]
    if ($aFmt[$r] -eq "General") {
        $ws.Cells.Item($r, 3).ClearFormats()
    } else {
        $ws.Cells.Item($r, 3).NumberFormat = $aFmt[$r]
    }
}

# --- Append new row 9 with the latest data point ---
$ws.Cells.Item(9, 1).Value2 = 806.651
$ws.Cells.Item(9, 2).Value2 = 823.9400000000001
$ws.Cells.Item(9, 3).Value2 = 45737

# Row 8's date (now in column C) reverts to the "regular" date style,
# while the newly appended row 9 takes on the "most recent" style that
# row 8 used to have.
$ws.Cells.Item(8, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 3).NumberFormat = "YYYY-MM-DD"
